$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numeric-looking strings such as '57.638.42',
# '6.00' or '0.997' that must stay literal text. Mark just those cells as
# Text-formatted first so Excel's COM Value setter does not coerce them
# into numbers (which would also collapse formatting like trailing zeros).
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '57.638.42'
$ws.Cells.Item(2, 5).Value = '  -4.16%  '
$ws.Cells.Item(3, 4).Value = '2.935.04'
$ws.Cells.Item(3, 5).Value = '  -2.10%  '
$ws.Cells.Item(4, 5).Value = '  -0.10%  '
$ws.Cells.Item(5, 4).Value = '549.68'
$ws.Cells.Item(5, 5).Value = '  -4.05%  '
$ws.Cells.Item(6, 4).Value = '130.81'
$ws.Cells.Item(6, 5).Value = '  +4.42%  '
$ws.Cells.Item(7, 5).Value = '  -0.04%  '
$ws.Cells.Item(8, 4).Value = '0.511'
$ws.Cells.Item(8, 5).Value = '  +1.88%  '
$ws.Cells.Item(9, 4).Value = '2.928.21'
$ws.Cells.Item(9, 5).Value = '  -2.18%  '
$ws.Cells.Item(10, 5).Value = '  -3.47%  '
$ws.Cells.Item(11, 4).Value = '4.78'
$ws.Cells.Item(11, 5).Value = '  -5.23%  '
$ws.Cells.Item(12, 5).Value = '  +1.52%  '
$ws.Cells.Item(13, 4).Value = '0.0000222'
$ws.Cells.Item(13, 5).Value = '  +0.78%  '
$ws.Cells.Item(14, 4).Value = '32.95'
$ws.Cells.Item(14, 5).Value = '  +1.47%  '
$ws.Cells.Item(15, 5).Value = '  +0.28%  '
$ws.Cells.Item(16, 4).Value = '3.418.04'
$ws.Cells.Item(16, 5).Value = '  -2.10%  '
$ws.Cells.Item(17, 4).Value = '6.86'
$ws.Cells.Item(17, 5).Value = '  +6.85%  '
$ws.Cells.Item(18, 4).Value = '2.929.45'
$ws.Cells.Item(18, 5).Value = '  -2.10%  '
$ws.Cells.Item(19, 4).Value = '57.595.95'
$ws.Cells.Item(19, 5).Value = '  -4.15%  '
$ws.Cells.Item(20, 4).Value = '418.27'
$ws.Cells.Item(20, 5).Value = '  -2.20%  '
$ws.Cells.Item(21, 4).Value = '13.20'
$ws.Cells.Item(21, 5).Value = '  +0.80%  '
$ws.Cells.Item(22, 4).Value = '0.687'
$ws.Cells.Item(22, 5).Value = '  +3.05%  '
$ws.Cells.Item(23, 4).Value = '6.99'
$ws.Cells.Item(23, 5).Value = '  -0.73%  '
$ws.Cells.Item(24, 4).Value = '13.06'
$ws.Cells.Item(24, 5).Value = '  +1.61%  '
$ws.Cells.Item(25, 5).Value = '  +0.77%  '
$ws.Cells.Item(26, 5).Value = '  -0.05%  '
$ws.Cells.Item(27, 4).Value = '0.997'
$ws.Cells.Item(27, 5).Value = '  -0.18%  '
$ws.Cells.Item(28, 5).Value = '  -2.54%  '
$ws.Cells.Item(29, 4).Value = '7.49'
$ws.Cells.Item(29, 5).Value = '  +3.48%  '
$ws.Cells.Item(30, 4).Value = '1.99'
$ws.Cells.Item(30, 5).Value = '  +2.04%  '
$ws.Cells.Item(31, 4).Value = '25.23'
$ws.Cells.Item(31, 5).Value = '  +0.12%  '
$ws.Cells.Item(32, 4).Value = '6.00'
$ws.Cells.Item(32, 5).Value = '  -2.39%  '
$ws.Cells.Item(33, 4).Value = '0.0971'
$ws.Cells.Item(33, 5).Value = '  +3.25%  '
$ws.Cells.Item(34, 4).Value = '5.67'
$ws.Cells.Item(34, 5).Value = '  +1.28%  '
$ws.Cells.Item(35, 4).Value = '0.941'
$ws.Cells.Item(35, 5).Value = '  +0.59%  '
$ws.Cells.Item(36, 4).Value = '2.08'
$ws.Cells.Item(36, 5).Value = '  +1.15%  '
$ws.Cells.Item(37, 2).Value = 'OKB'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(37, 4).Value = '47.97'
$ws.Cells.Item(37, 5).Value = '  -4.49%  '
$ws.Cells.Item(38, 2).Value = 'Cosmos'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(38, 4).Value = '8.74'
$ws.Cells.Item(38, 5).Value = '  +4.03%  '
$ws.Cells.Item(39, 4).Value = '0.0₃0683'
$ws.Cells.Item(39, 5).Value = '  +2.61%  '
$ws.Cells.Item(40, 4).Value = '2.56'
$ws.Cells.Item(40, 5).Value = '  +4.12%  '
$ws.Cells.Item(41, 5).Value = '  -0.22%  '
$ws.Cells.Item(42, 2).Value = 'Bittensor'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(42, 4).Value = '376.25'
$ws.Cells.Item(42, 5).Value = '  +0.15%  '
$ws.Cells.Item(43, 2).Value = 'VeChain'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(43, 4).Value = '0.0346'
$ws.Cells.Item(43, 5).Value = '  -2.53%  '
$ws.Cells.Item(44, 4).Value = '2.682.59'
$ws.Cells.Item(44, 5).Value = '  +0.18%  '
$ws.Cells.Item(45, 5).Value = '  +0.03%  '
$ws.Cells.Item(46, 5).Value = '  +1.72%  '
$ws.Cells.Item(47, 4).Value = '122.12'
$ws.Cells.Item(47, 5).Value = '  +1.31%  '
$ws.Cells.Item(48, 5).Value = '  +1.73%  '
$ws.Cells.Item(49, 4).Value = '1.98'
$ws.Cells.Item(49, 5).Value = '  -0.90%  '
$ws.Cells.Item(50, 4).Value = '23.21'
$ws.Cells.Item(50, 5).Value = '  -1.27%  '
$ws.Cells.Item(51, 5).Value = '  -0.03%  '
